$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowIndex, Prenom, Nom, P1..P10 (C..L)
$data = @(
    @(4,'Marie','Martin',1,1,1,1,1,1,1,1,0,1),
    @(5,'Jean','Bernard',1,0,0,1,1,0,0,0,1,1),
    @(6,'Sophie','Dubois',1,0,1,0,1,1,1,0,1,0),
    @(7,'Pierre','Thomas',0,1,1,1,0,1,1,1,1,1),
    @(8,'Julie','Robert',1,1,1,1,1,0,1,1,1,1),
    @(9,'Marc','Richard',0,1,1,1,0,1,0,1,1,0),
    @(10,'Laura','Petit',1,1,0,1,0,1,1,1,0,1),
    @(11,'Thomas','Durand',1,1,1,1,1,1,0,1,1,1),
    @(12,'Céline','Leroy',0,1,1,1,1,1,1,1,1,1),
    @(13,'Nicolas','Moreau',1,1,0,0,0,0,0,1,0,0),
    @(14,'Emma','Simon',0,1,0,1,1,0,1,1,0,1),
    @(15,'Lucas','Laurent',0,1,1,1,1,0,1,0,1,0),
    @(16,'Camille','Lefebvre',1,1,1,0,1,1,1,1,0,1),
    @(17,'Alexandre','Michel',1,1,1,1,1,1,1,1,1,1),
    @(18,'Léa','Garcia',1,1,0,1,1,1,0,1,1,1),
    @(19,'Julien','David',0,1,0,1,1,1,1,0,1,1),
    @(20,'Sarah','Bertrand',1,1,1,1,1,1,0,0,1,1),
    @(21,'Mathieu','Roux',1,1,1,0,1,1,1,1,1,1),
    @(22,'Chloé','Vincent',1,0,1,0,1,1,1,1,1,1),
    @(23,'David','Fournier',1,1,1,1,0,1,1,1,1,1),
    @(24,'Manon','Morel',1,0,1,1,1,1,1,1,0,1),
    @(25,'Antoine','Girard',1,0,0,1,1,1,1,1,0,1),
    @(26,'Charlotte','André',0,1,0,1,1,1,1,1,1,0),
    @(27,'Vincent','Lefevre',1,1,1,0,1,1,1,1,0,0),
    @(28,'Lisa','Mercier',0,0,0,1,1,1,1,1,0,0),
    @(29,'Maxime','Dupont',0,1,1,1,1,1,1,0,1,1),
    @(30,'Océane','Lambert',0,1,1,1,1,0,1,1,1,1),
    @(31,'Romain','Bonnet',1,1,0,0,0,1,0,1,1,1),
    @(32,'Eva','François',0,1,1,0,1,1,0,1,1,1),
    @(33,'Benjamin','Martinez',1,1,1,0,1,0,1,1,0,1),
    @(34,'Alice','Legrand',1,1,1,1,0,0,1,0,1,0),
    @(35,'François','Garnier',0,1,1,1,1,1,1,1,1,0),
    @(36,'Clara','Faure',0,1,0,1,1,1,0,1,1,1),
    @(37,'Hugo','Rousseau',1,1,1,0,1,1,1,1,0,1),
    @(38,'Inès','Blanc',1,1,1,1,1,1,1,0,0,1),
    @(39,'Arthur','Guerin',1,0,0,0,1,1,1,0,1,0),
    @(40,'Jade','Muller',1,1,0,1,1,1,1,1,0,0),
    @(41,'Louis','Henry',1,0,0,0,1,1,1,0,1,1),
    @(42,'Anaïs','Roussel',1,1,1,1,1,0,1,1,0,1),
    @(43,'Paul','Nicolas',0,1,1,1,0,1,1,0,0,0),
    @(44,'Lucie','Perrin',0,1,1,1,1,0,1,1,1,1),
    @(45,'Simon','Morin',1,0,1,1,1,1,1,1,1,0),
    @(46,'Margot','Mathieu',1,1,1,1,0,0,1,1,1,1),
    @(47,'Gabriel','Clement',0,1,1,0,0,0,1,1,1,1),
    @(48,'Zoé','Gauthier',0,1,0,1,0,1,0,1,0,1),
    @(49,'Raphaël','Dumont',1,1,0,1,1,1,1,1,0,0),
    @(50,'Louise','Lopez',0,1,1,1,1,1,1,1,0,0),
    @(51,'Tom','Fontaine',1,0,1,1,0,0,1,0,1,1),
    @(52,'Lina','Chevalier',0,0,1,1,1,1,1,1,0,0),
    @(53,'Nathan','Robin',1,1,0,1,0,1,1,1,1,1)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item($r, 3 + $i).Value = $entry[3 + $i]
    }
}
